$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change "True" values to "Yes" in B2:D2 and B3:D3 (leading apostrophe preserves
# the existing quote-prefixed text style these cells already use)
$ws.Range("B2").Formula = "'Yes"
$ws.Range("C2").Formula = "'Yes"
$ws.Range("D2").Formula = "'Yes"
$ws.Range("B3").Formula = "'Yes"
$ws.Range("C3").Formula = "'Yes"
$ws.Range("D3").Formula = "'Yes"

# Update the selected cell to reflect the new active cell (G6)
$ws.Range("G6").Select()
